# Adds three new data rows (124-126) to each of the four worksheets,
# matching the new daily log entries, and extends the used range
# (A1:I123 -> A1:I126) accordingly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("MID_LFT_#1")
$ws.Cells.Item(124, 1).Value = 45910.46265046296
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x01,0x90"
$ws.Cells.Item(124, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(124, 4).Value = "0x00,0xF0"
$ws.Cells.Item(124, 5).Value = "0x07"
$ws.Cells.Item(124, 6).Value = 400
$ws.Cells.Item(124, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(124, 8).Value = 240
$ws.Cells.Item(124, 9).Value = 7

$ws.Cells.Item(125, 1).Value = 45911.46655092593
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x01,0x90"
$ws.Cells.Item(125, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(125, 4).Value = "0x00,0xEC"
$ws.Cells.Item(125, 5).Value = "0x07"
$ws.Cells.Item(125, 6).Value = 400
$ws.Cells.Item(125, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(125, 8).Value = 240
$ws.Cells.Item(125, 9).Value = 7

$ws.Cells.Item(126, 1).Value = 45912.46122685185
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x01,0x90"
$ws.Cells.Item(126, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(126, 4).Value = "0x00,0xEC"
$ws.Cells.Item(126, 5).Value = "0x07"
$ws.Cells.Item(126, 6).Value = 400
$ws.Cells.Item(126, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(126, 8).Value = 236
$ws.Cells.Item(126, 9).Value = 7


$ws = $wb.Worksheets.Item("MID_LFT_#2")
$ws.Cells.Item(124, 1).Value = 45910.46265046296
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x01,0x7c"
$ws.Cells.Item(124, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(124, 4).Value = "0x01,0x04"
$ws.Cells.Item(124, 5).Value = "0x19"
$ws.Cells.Item(124, 6).Value = 380
$ws.Cells.Item(124, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(124, 8).Value = 260
$ws.Cells.Item(124, 9).Value = 25

$ws.Cells.Item(125, 1).Value = 45911.46655092593
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x01,0x7c"
$ws.Cells.Item(125, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(125, 4).Value = "0x01,0x00"
$ws.Cells.Item(125, 5).Value = "0x19"
$ws.Cells.Item(125, 6).Value = 380
$ws.Cells.Item(125, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(125, 8).Value = 256
$ws.Cells.Item(125, 9).Value = 25

$ws.Cells.Item(126, 1).Value = 45912.46122685185
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x01,0x7c"
$ws.Cells.Item(126, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(126, 4).Value = "0x01,0x00"
$ws.Cells.Item(126, 5).Value = "0x19"
$ws.Cells.Item(126, 6).Value = 380
$ws.Cells.Item(126, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(126, 8).Value = 256
$ws.Cells.Item(126, 9).Value = 25


$ws = $wb.Worksheets.Item("MID_PLT_#1")
$ws.Cells.Item(124, 1).Value = 45910.46265046296
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x00,0x6e"
$ws.Cells.Item(124, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(124, 4).Value = "0x00,0x58"
$ws.Cells.Item(124, 5).Value = "0x15"
$ws.Cells.Item(124, 6).Value = 110
$ws.Cells.Item(124, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(124, 8).Value = 88
$ws.Cells.Item(124, 9).Value = 15

$ws.Cells.Item(125, 1).Value = 45911.46655092593
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x00,0x6e"
$ws.Cells.Item(125, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(125, 4).Value = "0x00,0x58"
$ws.Cells.Item(125, 5).Value = "0x15"
$ws.Cells.Item(125, 6).Value = 110
$ws.Cells.Item(125, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(125, 8).Value = 88
$ws.Cells.Item(125, 9).Value = 15

$ws.Cells.Item(126, 1).Value = 45912.46122685185
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x00,0x6e"
$ws.Cells.Item(126, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(126, 4).Value = "0x00,0x58"
$ws.Cells.Item(126, 5).Value = "0x15"
$ws.Cells.Item(126, 6).Value = 110
$ws.Cells.Item(126, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(126, 8).Value = 88
$ws.Cells.Item(126, 9).Value = 15


$ws = $wb.Worksheets.Item("MID_PLT_#2")
$ws.Cells.Item(124, 1).Value = 45910.46265046296
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x00,0x82"
$ws.Cells.Item(124, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(124, 4).Value = "0x00,0x6F"
$ws.Cells.Item(124, 5).Value = "0x9"
$ws.Cells.Item(124, 6).Value = 130
$ws.Cells.Item(124, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(124, 8).Value = 111
$ws.Cells.Item(124, 9).Value = 9

$ws.Cells.Item(125, 1).Value = 45911.46655092593
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x00,0x82"
$ws.Cells.Item(125, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(125, 4).Value = "0x00,0x6F"
$ws.Cells.Item(125, 5).Value = "0x9"
$ws.Cells.Item(125, 6).Value = 130
$ws.Cells.Item(125, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(125, 8).Value = 111
$ws.Cells.Item(125, 9).Value = 9

$ws.Cells.Item(126, 1).Value = 45912.46122685185
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x00,0x82"
$ws.Cells.Item(126, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(126, 4).Value = "0x00,0x6E"
$ws.Cells.Item(126, 5).Value = "0x9"
$ws.Cells.Item(126, 6).Value = 130
$ws.Cells.Item(126, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(126, 8).Value = 110
$ws.Cells.Item(126, 9).Value = 9


Write-Output "Added rows 124-126 to all four sheets"
